# Facilitators guidelines - Surface Tension.docx
# Translate the English table labels/captions to Swahili (Kenya).
#
# wdReplaceAll = 2 is used throughout so that every matching occurrence
# in the document body gets updated in a single call; the one phrase
# that repeats ("Assist the process, provoke thoughts") is therefore
# handled without any extra bookkeeping.
#
# NOTE: "General VMC Video Introduction" contains "Video Introduction"
# as a substring, so it is translated *before* the shorter "Video
# Introduction" string to avoid a partial/incorrect match.

$d = $word.ActiveDocument

function Translate-All($old, $new) {
    $d.Content.Find.Execute(
        $old,   # FindText
        $true,  # MatchCase
        $false, # MatchWholeWord
        $false, # MatchWildcards
        $false, # MatchSoundsLike
        $false, # MatchAllWordForms
        $true,  # Forward
        1,      # Wrap (wdFindContinue)
        $false, # Format
        $new,   # ReplaceWith
        2       # Replace (wdReplaceAll)
    ) | Out-Null
}

Translate-All "Video Title" "Kichwa cha Video"
Translate-All "Topic" "Mada"
Translate-All "Aim(s)" "Malengo"
Translate-All "Length" "Urefu"
Translate-All "Camp Location" "Mahali pa Kambi"
Translate-All "Facilitators" "Wawezeshaji"
Translate-All "N. of students" "N. ya wanafunzi"
Translate-All "Date" "Tarehe"
Translate-All "Resources" "Rasilimali"
Translate-All "needed" "inahitajika"
Translate-All "Preparations" "Maandalizi"
Translate-All "Video time" "Muda wa video"
Translate-All "What facilitator does" "Mwezeshaji anafanya nini"
Translate-All "What learners do" "Wanachofanya wanafunzi"
Translate-All "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Translate-All "Video Introduction" "Utangulizi wa Video"
Translate-All "Introduction of the first experiment" "Utangulizi wa jaribio la kwanza"
Translate-All "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
